# Add the 2023 row to both "Hoja1" and "Hoja2" sheets, matching the
# pattern already used for the previous years (file_path in column A,
# file_name/year in column B).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Hoja1
$ws2 = $wb.Worksheets.Item(2)   # Hoja2

$filePath = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\2023.xlsx"
$fileYear = "2023"

foreach ($ws in @($ws1, $ws2)) {
    # Column B needs to stay text (e.g. "2023", not the number 2023) just
    # like the rest of the column, so force the text format before writing.
    $ws.Range("B23").NumberFormat = "@"

    $ws.Range("A23").Value = $filePath
    $ws.Range("B23").Value = $fileYear
}

# Update the selections to cover the new last row, and re-activate Hoja2
# then Hoja1 last so Hoja1 ends up the selected tab again (matching the
# original workbook state).
$ws2.Activate()
$ws2.Range("A2:B23").Select()

$ws1.Activate()
$ws1.Range("A2:B23").Select()

# Match the resized/repositioned workbook window from the saved file.
$win = $wb.Windows.Item(1)
$win.Left = 1950
$win.Top = 1950
$win.Width = 18075
$win.Height = 12255
